$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.575.93"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.667.83"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.84"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.80"
$ws.Range("E6").Value = "  -1.17%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.625"
$ws.Range("E8").Value = "  +5.80%  "

$ws.Range("E9").Value = "  +2.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.401"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("E11").Value = "  -3.95%  "

$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.33"
$ws.Range("E13").Value = "  -3.92%  "

$ws.Range("E14").Value = "  -3.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.145.67"
$ws.Range("E15").Value = "  -1.82%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.447.58"
$ws.Range("E16").Value = "  -0.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.644.90"
$ws.Range("E17").Value = "  -3.15%  "

$ws.Range("E18").Value = "  -0.34%  "

$ws.Range("E19").Value = "  -2.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("E20").Value = "  -1.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.72"
$ws.Range("E21").Value = "  -3.39%  "

$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.56"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000112"
$ws.Range("E24").Value = "  +3.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.66"
$ws.Range("E25").Value = "  -1.74%  "

$ws.Range("E26").Value = "  -1.78%  "

$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value = "  -7.00%  "

$ws.Range("E29").Value = "  -2.50%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.15"
$ws.Range("E31").Value = "  -2.73%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "530.61"
$ws.Range("E32").Value = "  -2.10%  "

$ws.Range("E33").Value = "  -3.36%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.46"
$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.423"
$ws.Range("E36").Value = "  -3.24%  "

$ws.Range("E37").Value = "  -2.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "158.32"
$ws.Range("E39").Value = "  -3.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  -3.69%  "

$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "164.39"
$ws.Range("E42").Value = "  -4.07%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.14"
$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +0.50%  "

$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0607"
$ws.Range("E45").Value = "  -1.65%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.81"
$ws.Range("E46").Value = "  -3.56%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.640"
$ws.Range("E47").Value = "  -3.40%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0258"
$ws.Range("E48").Value = "  -3.34%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0259"
$ws.Range("E49").Value = "  +12.09%  "

$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.01"
$ws.Range("E51").Value = "  -4.80%  "
